# Added check for empty Then blocks.
#
# The "Workflow" sheet's row 21 used to describe a generic, slightly
# mis-aligned "Empty Then block" check (values were shifted by one
# shared-string relative to every other row because of an authoring
# mistake). This edit rewrites row 21 with the corrected / expanded
# wording (mentioning FlowDecision's True branch as well) and moves the
# on-screen selection down to that row, the way Excel would leave it
# after the author finished editing there.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workflow")

# --- Row 21: "Empty Then block" check -------------------------------------
$ws.Range("A21").Value = "Yes"
$ws.Range("B21").Value = "Empty Then block"
$ws.Range("C21").Value = "Checks\EmptyThenBlock.xaml"
$ws.Range("F21").Value = "In an If activity (or FlowDecision activity), it is recommended to define the condition so that the Then block (True branch) is always used. Having an empty Then block (True branch) makes the workflow harder to understand."
$ws.Range("G21").Value = "Redefine the condition so that actions are taken in the Then block (True branch, in case of FlowDecision)."

# The longer explanation text now wraps onto one more line, so the
# (wrap-text, auto-height) row grows from 75 to 93.75 points.
$ws.Rows.Item(21).RowHeight = 93.75

# Leave the selection on B21, scrolled so row 21 is in view, matching
# where the author was last working.
$ws.Activate()
$ws.Range("B21").Select()
